$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Convert the "Correct Answer" (H) and "Time in seconds" (I) columns
# from text values to real numbers for rows 2-9.
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 30

$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 45

$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 45

$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 30

$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 30

$ws.Range("H7").Value = 2
$ws.Range("I7").Value = 45

$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 30

$ws.Range("H9").Value = 2
$ws.Range("I9").Value = 30

# Row 10: update the fill-in-the-blank answer text, clear the now-unused
# "Correct Answer" cell, and make "Time in seconds" a real number.
# Force text storage so "1.0" isn't auto-coerced into the number 1.
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "1.0"
$ws.Range("H10").Value = ""
$ws.Range("I10").Value = 45

# Row 11: "Time in seconds" becomes a real number.
$ws.Range("I11").Value = 60
